# This edit reorders the data rows 2-7 of the "Artfynd" sheet.
# Each whole row's content (all its populated columns) moves together to a
# new row position; row 1 (header) and rows 8-14 are untouched.
#
# Resulting permutation (after-row -> before-row the data came from):
#   2 <- 4
#   3 <- 7
#   4 <- 2
#   5 <- 3
#   6 <- 5
#   7 <- 6
#
# i.e. two cycles: (2 4) and (3 7 6 5)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Contiguous column runs that are actually populated for the data rows
# (keeps us from creating stray empty <c> cells in unused columns such as
# K-O, X, AC, AF, AH, AJ-AM, AP-AS, AU, AV).
$ColumnRuns = @(
    @("A", "J"),
    @("P", "W"),
    @("Y", "AB"),
    @("AD", "AE"),
    @("AG", "AG"),
    @("AN", "AO"),
    @("AT", "AT"),
    @("AW", "AY")
)

function Copy-RowData {
    param([int]$SrcRow, [int]$DstRow)

    foreach ($run in $ColumnRuns) {
        $c1 = $run[0]
        $c2 = $run[1]
        $ws.Range("${c1}${SrcRow}:${c2}${SrcRow}").Copy($ws.Range("${c1}${DstRow}:${c2}${DstRow}"))
    }

    # Column AI ("Biotop") is sparse - only some rows have it. Always clear
    # the destination first, then copy it across only if the source has it.
    $ws.Range("AI${DstRow}").ClearContents()
    $aiVal = $ws.Range("AI${SrcRow}").Value()
    if ($aiVal -ne $null -and $aiVal -ne "") {
        $ws.Range("AI${SrcRow}").Copy($ws.Range("AI${DstRow}"))
    }
}

function Clear-RowData {
    param([int]$Row)

    foreach ($run in $ColumnRuns) {
        $c1 = $run[0]
        $c2 = $run[1]
        $ws.Range("${c1}${Row}:${c2}${Row}").ClearContents()
    }
    $ws.Range("AI${Row}").ClearContents()
}

# Scratch rows far away from any real data, used to hold a row's content
# temporarily while permuting.
$Scratch1 = 1000
$Scratch2 = 1001

# Cycle (2 4): new-row2 = old-row4, new-row4 = old-row2
Copy-RowData 2 $Scratch1
Clear-RowData 2
Copy-RowData 4 2
Clear-RowData 4
Copy-RowData $Scratch1 4
Clear-RowData $Scratch1

# Cycle (3 7 6 5): new-row3 = old-row7, new-row7 = old-row6,
#                  new-row6 = old-row5, new-row5 = old-row3
Copy-RowData 3 $Scratch2
Clear-RowData 3
Copy-RowData 7 3
Clear-RowData 7
Copy-RowData 6 7
Clear-RowData 6
Copy-RowData 5 6
Clear-RowData 5
Copy-RowData $Scratch2 5
Clear-RowData $Scratch2
